$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old rows 5-7 entirely (they no longer exist in the new data range)
$ws.Rows("5:7").Delete()

# New data for rows 1-4
$data = @(
    @(2, "pan", 44819),
    @(3, "Huevos", 44819),
    @(4, "Arroz", 44819),
    @(5, "Mantequilla", 44819)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 3).NumberFormat = "yyyy-mm-dd"
}
